$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "0.450", "0.999") keep their exact textual representation instead
# of being coerced into numbers (which would drop trailing zeros, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.743.02"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "3.815.67"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "601.27"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "165.90"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "6.29"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").Value = "0.450"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  -0.94%  "
$ws.Range("D13").Value = "35.75"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "4.455.80"
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").Value = "3.811.61"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").Value = "67.770.04"
$ws.Range("E16").Value = "  +0.45%  "
$ws.Range("D17").Value = "18.39"
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "463.12"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  -4.17%  "
$ws.Range("D24").Value = "83.38"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").Value = "12.06"
$ws.Range("E25").Value = "  +1.52%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "3.964.94"
$ws.Range("E29").Value = "  +1.67%  "
$ws.Range("D30").Value = "2.78"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "7.36"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("D33").Value = "29.53"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "9.07"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").Value = "0.0998"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "0.138"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.27"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "0.997"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("D40").Value = "5.79"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("D43").Value = "48.08"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "0.299"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("E45").Value = "  +9.34%  "
$ws.Range("D46").Value = "43.11"
$ws.Range("E46").Value = "  -4.82%  "
$ws.Range("E47").Value = "  +12.58%  "
$ws.Range("D48").Value = "148.86"
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("D51").Value = "385.82"
$ws.Range("E51").Value = "  -1.10%  "
